# Auto-generated edit script: rebuilds the "About Tati" body content
$d = $word.ActiveDocument

$paraXml = @(
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+RiyBIaSwgSeKAmW0gKipUYXRpYW5hIEJyaW1tKiosIGEgMjQteWVhci1vbGQgQ2xpbmljYWwgSW5mb3JtYXRpY3MgbWFqb3IgYW5kIHBhc3Npb25hdGUgdGVjaG5vbG9naXN0Ljwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+OkyBJ4oCZbSBjdXJyZW50bHkgcHVyc3VpbmcgbXkgQmFjaGVsb3LigJlzIGRlZ3JlZSBhbmQgcGxhbiB0byBlYXJuIG15IE1hc3RlcuKAmXMgaW4gQ29tcHV0ZXIgU2NpZW5jZS4gTXkgZ29hbCBpcyB0byBiZWNvbWUgYSAqKmZ1bGwtc3RhY2sgZGV2ZWxvcGVyKiogc3BlY2lhbGl6aW5nIGluOjwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBBSSBkZXZlbG9wbWVudCAmYW1wOyBhdXRvbWF0aW9uPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBHYW1lIGRldmVsb3BtZW50IChlbXVsYXRvcnMsIHJldHJvIHN5c3RlbXMpPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBEYXRhIGFuYWx5dGljcyAmYW1wOyBtYWNoaW5lIGxlYXJuaW5nPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBDdXN0b20gZGlnaXRhbCB0b29scyBmb3Igc21hbGwgJmFtcDsgbGFyZ2UgYnVzaW5lc3Nlczwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+boO+4jyBJIGJ1aWx0IHRoaXMgZW50aXJlIHBvcnRmb2xpbyBzaXRlIGZyb20gc2NyYXRjaCB1c2luZyAqKlZpc3VhbCBTdHVkaW8gQ29kZSAoVlMgQ29kZSkqKi4gSXTigJlzIGFuIG9uZ29pbmcgcHJvamVjdCBhbmQgd2lsbCBiZSB1cGRhdGVkIGZyZXF1ZW50bHkgd2l0aCBuZXcgcHJvamVjdHMsIFVJIGVuaGFuY2VtZW50cywgYW5kIGZlYXR1cmVzLjwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+SvCAqKlNlcnZpY2VzIENvbWluZyBGYWxsIDIwMjUqKjo8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+VGF0aWFuYSBpcyBjdXJyZW50bHkgYnVpbGRpbmcgYSBzdWl0ZSBvZiBkaWdpdGFsIHByb2R1Y3RzIHRoYXQgd2lsbCBiZSBhdmFpbGFibGUgZm9yIHB1cmNoYXNlLiBUaGVzZSBpbmNsdWRlOjwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBBSS1wb3dlcmVkIHdvcmtmbG93IGF1dG9tYXRpb25zPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBDdXN0b20gQUkgYWdlbnRzPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBNaW5pLXN5c3RlbXMgd2l0aCBmcm9udGVuZCBVSXM8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBEZXZlbG9wZXIgdG9vbHMgdG8gaGVscCBzbWFsbCBidXNpbmVzc2VzIGdyb3cgYW5kIHNjYWxlPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+MjSBJIGV2ZW50dWFsbHkgcGxhbiB0bzo8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBCdWlsZCBmcmVlIEFJIGFuZCB0ZWNoIGNvdXJzZXMgZm9yIHVuZGVyLXJlc291cmNlZCBjb21tdW5pdGllczwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBTdXBwb3J0IG5vbnByb2ZpdHMgZm9jdXNlZCBvbiBhbmltYWxzLCBlZHVjYXRpb24sIGFuZCB0aGUgZW52aXJvbm1lbnQ8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBDcmVhdGUgYSBzYW5jdHVhcnkgd2hlcmUgc2hlbHRlciBhbmltYWxzIGNhbiBsaXZlIHRoZWlyIGZpbmFsIGRheXMgaW4gcGVhY2U8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+OuCBGdW4gZmFjdHMgYWJvdXQgbWU6PC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBJ4oCZbSB0aGUgb2xkZXN0IGdpcmwgb2YgNiBzaWJsaW5nczwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBJIGxvdmUgYW5pbWFscyBhbmQgd2FudCB0byByZXNjdWUgYXMgbWFueSBhcyBJIGNhbjwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBJ4oCZbSBsZWFybmluZyB0byBwbGF5IGd1aXRhciAoYW5kIGRydW1zIG5leHQhKTwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBJIHdvcmsgb3V0IHNvbWV0aW1lcywgYnV0IEnigJltIGEgKipodWdlIGZvb2RpZSoqPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBNeSBmYW1pbHkgaXMgSmFtYWljYW4g8J+Hr/Cfh7I8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBNeSBmYXZvcml0ZSBib29rIGlzICpUaGUgU2V2ZW4gSHVzYmFuZHMgb2YgRXZlbHluIEh1Z28qPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+4pyoIEnigJltIGVhc3lnb2luZywgY3VyaW91cywgYW5kIGRlZXBseSBwYXNzaW9uYXRlIGFib3V0IHdoYXQgSSBkby4gV2hldGhlciBpdOKAmXMgYnVpbGRpbmcgYXV0b21hdGlvbiBzeXN0ZW1zLCBoZWxwaW5nIHBlb3BsZSBnZXQgdGVjaC1zYXZ2eSwgb3IgYnJpbmdpbmcgYmlnIGlkZWFzIHRvIGxpZmUg4oCUIEnigJltIGhlcmUgdG8gY3JlYXRlIHNvbWV0aGluZyBtZWFuaW5nZnVsIGFuZCBmdW4uPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LS0tPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+8J+kliAqKjwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxTdGFydCIvPjx3OnI+PHc6clByPjx3OnJGb250cyB3OmFzY2lpPSJTZWdvZSBVSSBFbW9qaSIgdzpoQW5zaT0iU2Vnb2UgVUkgRW1vamkiIHc6Y3M9IlNlZ29lIFVJIEVtb2ppIi8+PC93OnJQcj48dzp0PkJyaW1tQm90PC93OnQ+PC93OnI+PHc6cHJvb2ZFcnIgdzp0eXBlPSJzcGVsbEVuZCIvPjx3OnI+PHc6clByPjx3OnJGb250cyB3OmFzY2lpPSJTZWdvZSBVSSBFbW9qaSIgdzpoQW5zaT0iU2Vnb2UgVUkgRW1vamkiIHc6Y3M9IlNlZ29lIFVJIEVtb2ppIi8+PC93OnJQcj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiBJbnN0cnVjdGlvbmFsIEJlaGF2aW9yKio8L3c6dD48L3c6cj48L3c6cD4=",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBJZiBzb21lb25lIGFza3MgdW5yZWxhdGVkLCBpbmFwcHJvcHJpYXRlLCBvciDigJx3ZWlyZOKAnSBxdWVzdGlvbnMsIGtpbmRseSByZWRpcmVjdCB0aGUgY29udmVyc2F0aW9uOjwvdzp0PjwvdzpyPjwvdzpwPg==",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+ICAmZ3Q7IOKAnEknbSA8L3c6dD48L3c6cj48dzpwcm9vZkVyciB3OnR5cGU9InNwZWxsU3RhcnQiLz48dzpyPjx3OnJQcj48dzpyRm9udHMgdzphc2NpaT0iU2Vnb2UgVUkgRW1vamkiIHc6aEFuc2k9IlNlZ29lIFVJIEVtb2ppIiB3OmNzPSJTZWdvZSBVSSBFbW9qaSIvPjwvdzpyUHI+PHc6dD5CcmltbUJvdDwvdzp0PjwvdzpyPjx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxFbmQiLz48dzpyPjx3OnJQcj48dzpyRm9udHMgdzphc2NpaT0iU2Vnb2UgVUkgRW1vamkiIHc6aEFuc2k9IlNlZ29lIFVJIEVtb2ppIiB3OmNzPSJTZWdvZSBVSSBFbW9qaSIvPjwvdzpyUHI+PHc6dD4sIGhlcmUgdG8gYW5zd2VyIHF1ZXN0aW9ucyBhYm91dCBUYXRpYW5hLCBoZXIgd29yaywgYW5kIHRoaXMgcG9ydGZvbGlvISDwn5iKIEZlZWwgZnJlZSB0byBhc2sgbWUgYWJvdXQgaGVyIGJhY2tncm91bmQsIHByb2plY3RzLCBvciBmdXR1cmUgcGxhbnMu4oCdPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBQcmlvcml0aXplIHNoYXJpbmcgdXNlZnVsIGluZm8gYWJvdXQgVGF0aWFuYeKAmXMgc2tpbGxzLCBpbnRlcmVzdHMsIHByb2plY3RzLCBhbmQgc2VydmljZXMuPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3OnQ+LSBTdGF5IGZyaWVuZGx5LCBjaGlsbCwgYW5kIGNvbnZlcnNhdGlvbmFsIOKAlCBsaWtlIFRhdGlhbmEhPC93OnQ+PC93OnI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnBQcj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjwvdzpwUHI+PC93OnA+",
  "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6clByPjx3OnJGb250cyB3OmFzY2lpPSJTZWdvZSBVSSBFbW9qaSIgdzpoQW5zaT0iU2Vnb2UgVUkgRW1vamkiIHc6Y3M9IlNlZ29lIFVJIEVtb2ppIi8+PC93OnJQcj48dzp0PvCfk6wgV2FudCB0byBjb25uZWN0IG9yIGNvbGxhYm9yYXRlPyBVc2UgdGhlICoqY29udGFjdCBwYWdlKiogdG8gc2VuZCBhIG1lc3NhZ2UgZGlyZWN0bHkgdG8gVGF0aWFuYS48L3c6dD48L3c6cj48dzpyPjx3OnJQcj48dzpyRm9udHMgdzphc2NpaT0iU2Vnb2UgVUkgRW1vamkiIHc6aEFuc2k9IlNlZ29lIFVJIEVtb2ppIiB3OmNzPSJTZWdvZSBVSSBFbW9qaSIvPjwvdzpyUHI+PHc6YnIvPjwvdzpyPjx3OnI+PHc6clByPjx3OnJGb250cyB3OmFzY2lpPSJTZWdvZSBVSSBFbW9qaSIgdzpoQW5zaT0iU2Vnb2UgVUkgRW1vamkiIHc6Y3M9IlNlZ29lIFVJIEVtb2ppIi8+PC93OnJQcj48dzpici8+PC93OnI+PHc6cj48dzpyUHI+PHc6ckZvbnRzIHc6YXNjaWk9IlNlZ29lIFVJIEVtb2ppIiB3OmhBbnNpPSJTZWdvZSBVSSBFbW9qaSIgdzpjcz0iU2Vnb2UgVUkgRW1vamkiLz48L3c6clByPjx3Omxhc3RSZW5kZXJlZFBhZ2VCcmVhay8+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5XYW50IHRvIHNlZSBtb3JlIG9mIGhlciBwcm9qZWN0cz8gR28gdG8gdGhlIFByb2plY3RzIHBhZ2UgdG8gdmlldyBoZXIgY3VycmVudCBwcm9qZWN0cy4gU2hlIHdpbGwgYmUgdXBkYXRpbmcgdGhhdCBwYWdlIGZyZXF1ZW50bHkuIDwvdzp0PjwvdzpyPjwvdzpwPg=="
)

# 1) Collapse the document down to a single, empty paragraph.
#    (Word never allows deleting the very last paragraph mark, so we
#    delete everything up to the start of the last paragraph, then trim
#    the last paragraph's own text, leaving one empty <w:p/>.)
$n = $d.Paragraphs.Count
$lastStart = $d.Paragraphs.Item($n).Range.Start
if ($lastStart -gt 0) {
    $d.Range(0, $lastStart).Delete()
}
$p1 = $d.Paragraphs.Item(1)
if ($p1.Range.End - 1 -gt $p1.Range.Start) {
    $d.Range($p1.Range.Start, $p1.Range.End - 1).Delete()
}

# 2) Replace that single empty paragraph, then append one new empty
#    paragraph per remaining item and fill each via InsertXML so we can
#    control run/proofErr/break structure exactly.
$first = $true
foreach ($b64 in $paraXml) {
    $bytes = [Convert]::FromBase64String($b64)
    $xml = [System.Text.Encoding]::UTF8.GetString($bytes)
    if ($first) {
        $target = $d.Paragraphs.Item($d.Paragraphs.Count)
        $first = $false
    } else {
        $lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
        $lastP.Range.InsertParagraphAfter()
        $target = $d.Paragraphs.Item($d.Paragraphs.Count)
    }
    $target.Range.InsertXML($xml)
}

Write-Output "Done. Paragraphs=$($d.Paragraphs.Count)"
